$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Alex Schell" row (row 18) entirely, shifting the rows below it up.
$ws.Rows("18:18").Delete()

# Rename the "Shirt Size" header in B1 to "Polo Shirt Size"
$ws.Range("B1").Value = "Polo Shirt Size"

# Update selection to B1 (matches the post-edit active cell in the file)
$ws.Range("B1").Select()
